$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamps for rows 2 and 3
$ws.Range("A2").Value = "05/26/2021 01:55:32"
$ws.Range("A3").Value = "05/26/2021 01:56:34"

# Delete rows 4 through 24 (the trailing data rows)
$ws.Range("A4:B24").EntireRow.Delete()
